$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.323874002307434
$ws.Cells.Item(2, 4).Value = 0.1573910469263495
$ws.Cells.Item(2, 5).Value = 0.2106255498562959
$ws.Cells.Item(2, 6).Value = 1.848000169810419
$ws.Cells.Item(2, 7).Value = 1.341898598714067
$ws.Cells.Item(2, 8).Value = 1.220470549407452
$ws.Cells.Item(2, 9).Value = 1.071364875308631
$ws.Cells.Item(2, 10).Value = 0.3078422427925744
$ws.Cells.Item(2, 12).Value = 0.6734733331329892
$ws.Cells.Item(2, 13).Value = 0.4667821546985351

$ws.Cells.Item(3, 2).Value = 1.235439004134491
$ws.Cells.Item(3, 4).Value = 0.150591120474445
$ws.Cells.Item(3, 5).Value = 0.199351563426525
$ws.Cells.Item(3, 6).Value = 1.855328665068868
$ws.Cells.Item(3, 7).Value = 1.330238678936553
$ws.Cells.Item(3, 8).Value = 1.222864123586305
$ws.Cells.Item(3, 9).Value = 1.093660492588789
$ws.Cells.Item(3, 10).Value = 0.2895890715531664
$ws.Cells.Item(3, 12).Value = 0.6220564915655302
$ws.Cells.Item(3, 13).Value = 0.4335791087960317

$ws.Cells.Item(4, 2).Value = 1.181256805577306
$ws.Cells.Item(4, 4).Value = 0.1463762148200374
$ws.Cells.Item(4, 5).Value = 0.1923918558679532
$ws.Cells.Item(4, 6).Value = 1.861297682897373
$ws.Cells.Item(4, 7).Value = 1.324344595726032
$ws.Cells.Item(4, 8).Value = 1.225210864891849
$ws.Cells.Item(4, 9).Value = 1.108137494177141
$ws.Cells.Item(4, 10).Value = 0.2783374591724197
$ws.Cells.Item(4, 12).Value = 0.590577888846326
$ws.Cells.Item(4, 13).Value = 0.4132440228912131

$ws.Cells.Item(5, 2).Value = 1.159207787650644
$ws.Cells.Item(5, 4).Value = 0.1446485558091979
$ws.Cells.Item(5, 5).Value = 0.1895463900453578
$ws.Cells.Item(5, 6).Value = 1.864098394403527
$ws.Cells.Item(5, 7).Value = 1.322259166900892
$ws.Cells.Item(5, 8).Value = 1.226387114487537
$ws.Cells.Item(5, 9).Value = 1.114234744015439
$ws.Cells.Item(5, 10).Value = 0.2737414775024973
$ws.Cells.Item(5, 12).Value = 0.5777735137383218
$ws.Cells.Item(5, 13).Value = 0.4049706831977744

$ws.Cells.Item(6, 2).Value = 1.155548452682325
$ws.Cells.Item(6, 4).Value = 0.1443610705286034
$ws.Cells.Item(6, 5).Value = 0.1890733417507917
$ws.Cells.Item(6, 6).Value = 1.864585659762383
$ws.Cells.Item(6, 7).Value = 1.321931947408899
$ws.Cells.Item(6, 8).Value = 1.226595695366825
$ws.Cells.Item(6, 9).Value = 1.115259124470792
$ws.Cells.Item(6, 10).Value = 0.2729776692742121
$ws.Cells.Item(6, 12).Value = 0.5756487828058141
$ws.Cells.Item(6, 13).Value = 0.4035977201558794

$ws.Cells.Item(7, 2).Value = 1.180959319055574
$ws.Cells.Item(7, 4).Value = 0.1463529557795198
$ws.Cells.Item(7, 5).Value = 0.1923535185747696
$ws.Cells.Item(7, 6).Value = 1.861333964659607
$ws.Cells.Item(7, 7).Value = 1.324315191801318
$ws.Cells.Item(7, 8).Value = 1.225225838528388
$ws.Cells.Item(7, 9).Value = 1.108218923451497
$ws.Cells.Item(7, 10).Value = 0.2782755198111886
$ws.Cells.Item(7, 12).Value = 0.5904051092470581
$ws.Cells.Item(7, 13).Value = 0.4131323910996514

$ws.Cells.Item(8, 2).Value = 1.293357955405327
$ws.Cells.Item(8, 4).Value = 0.1550546373855184
$ws.Cells.Item(8, 5).Value = 0.2067460955419023
$ws.Cells.Item(8, 6).Value = 1.850221369799698
$ws.Cells.Item(8, 7).Value = 1.337614657763453
$ws.Cells.Item(8, 8).Value = 1.221113442317147
$ws.Cells.Item(8, 9).Value = 1.078888833289986
$ws.Cells.Item(8, 10).Value = 0.3015578309370852
$ws.Cells.Item(8, 12).Value = 0.6557260199969335
$ws.Cells.Item(8, 13).Value = 0.4553232177666544

$ws.Cells.Item(9, 2).Value = 1.514661959782757
$ws.Cells.Item(9, 4).Value = 0.1718069988422144
$ws.Cells.Item(9, 5).Value = 0.2346705909023328
$ws.Cells.Item(9, 6).Value = 1.84014299644403
$ws.Cells.Item(9, 7).Value = 1.373812565224284
$ws.Cells.Item(9, 8).Value = 1.22003762924794
$ws.Cells.Item(9, 9).Value = 1.027630535289772
$ws.Cells.Item(9, 10).Value = 0.3468568449776939
$ws.Cells.Item(9, 12).Value = 0.7845366167092607
$ws.Cells.Item(9, 13).Value = 0.5384581607820849

$ws.Cells.Item(10, 2).Value = 1.677761865392029
$ws.Cells.Item(10, 4).Value = 0.1839310722159269
$ws.Cells.Item(10, 5).Value = 0.2550033194935608
$ws.Cells.Item(10, 6).Value = 1.839958738120188
$ws.Cells.Item(10, 7).Value = 1.406690244815763
$ws.Cells.Item(10, 8).Value = 1.223550705037496
$ws.Cells.Item(10, 9).Value = 0.9937990857377947
$ws.Cells.Item(10, 10).Value = 0.379913251899751
$ws.Cells.Item(10, 12).Value = 0.8796080617011341
$ws.Cells.Item(10, 13).Value = 0.5997719802066968

$ws.Cells.Item(11, 2).Value = 1.752064614082883
$ws.Cells.Item(11, 4).Value = 0.1894080953029089
$ws.Cells.Item(11, 5).Value = 0.2642133602873074
$ws.Cells.Item(11, 6).Value = 1.841460185265447
$ws.Cells.Item(11, 7).Value = 1.423036699857391
$ws.Cells.Item(11, 8).Value = 1.226092869790023
$ws.Cells.Item(11, 9).Value = 0.9792423430363506
$ws.Cells.Item(11, 10).Value = 0.3949015090186521
$ws.Cells.Item(11, 12).Value = 0.9229532113089078
$ws.Cells.Item(11, 13).Value = 0.6277149180869941

$ws.Cells.Item(12, 2).Value = 1.780215794938783
$ws.Cells.Item(12, 4).Value = 0.1914766880261567
$ws.Cells.Item(12, 5).Value = 0.2676952511999531
$ws.Cells.Item(12, 6).Value = 1.842258087012098
$ws.Cells.Item(12, 7).Value = 1.429428551995784
$ws.Cells.Item(12, 8).Value = 1.227192019380823
$ws.Cells.Item(12, 9).Value = 0.9738502025044475
$ws.Cells.Item(12, 10).Value = 0.4005699246228573
$ws.Cells.Item(12, 12).Value = 0.9393806188907661
$ws.Cells.Item(12, 13).Value = 0.6383032766792951

$ws.Cells.Item(13, 2).Value = 1.774152314692685
$ws.Cells.Item(13, 4).Value = 0.1910314202946211
$ws.Cells.Item(13, 5).Value = 0.2669456207922138
$ws.Cells.Item(13, 6).Value = 1.842076022773796
$ws.Cells.Item(13, 7).Value = 1.428042946736781
$ws.Cells.Item(13, 8).Value = 1.226949215793638
$ws.Cells.Item(13, 9).Value = 0.9750061457402373
$ws.Cells.Item(13, 10).Value = 0.399349458672134
$ws.Cells.Item(13, 12).Value = 0.9358420861321406
$ws.Cells.Item(13, 13).Value = 0.6360225795726109

$ws.Cells.Item(14, 2).Value = 1.754380347758456
$ws.Cells.Item(14, 4).Value = 0.1895783882449962
$ws.Cells.Item(14, 5).Value = 0.2644999332853715
$ws.Cells.Item(14, 6).Value = 1.841521224962719
$ws.Cells.Item(14, 7).Value = 1.423558506048749
$ws.Cells.Item(14, 8).Value = 1.226180557859664
$ws.Cells.Item(14, 9).Value = 0.9787963177589871
$ws.Cells.Item(14, 10).Value = 0.395368000810322
$ws.Cells.Item(14, 12).Value = 0.924304434667846
$ws.Cells.Item(14, 13).Value = 0.6285858909443789

$ws.Cells.Item(15, 2).Value = 1.742271281057413
$ws.Cells.Item(15, 4).Value = 0.1886876581253887
$ws.Cells.Item(15, 5).Value = 0.2630011276181463
$ws.Cells.Item(15, 6).Value = 1.841211302658451
$ws.Cells.Item(15, 7).Value = 1.420837993837239
$ws.Cells.Item(15, 8).Value = 1.225727529206182
$ws.Cells.Item(15, 9).Value = 0.9811335691829415
$ws.Cells.Item(15, 10).Value = 0.392928284003716
$ws.Cells.Item(15, 12).Value = 0.9172390428685162
$ws.Cells.Item(15, 13).Value = 0.6240316023239814

$ws.Cells.Item(16, 2).Value = 1.672908075412408
$ws.Cells.Item(16, 4).Value = 0.1835723713313655
$ws.Cells.Item(16, 5).Value = 0.254400621268303
$ws.Cells.Item(16, 6).Value = 1.839892632129178
$ws.Cells.Item(16, 7).Value = 1.40565010192617
$ws.Cells.Item(16, 8).Value = 1.223403627852207
$ws.Cells.Item(16, 9).Value = 0.9947672018375515
$ws.Cells.Item(16, 10).Value = 0.3789327230857964
$ws.Cells.Item(16, 12).Value = 0.8767772683046928
$ws.Cells.Item(16, 13).Value = 0.5979468386125291

$ws.Cells.Item(17, 2).Value = 1.630382769838832
$ws.Cells.Item(17, 4).Value = 0.1804245518699616
$ws.Cells.Item(17, 5).Value = 0.249114332791379
$ws.Cells.Item(17, 6).Value = 1.839490672121357
$ws.Cells.Item(17, 7).Value = 1.396690235766982
$ws.Cells.Item(17, 8).Value = 1.222220303268443
$ws.Cells.Item(17, 9).Value = 1.003344680798897
$ws.Cells.Item(17, 10).Value = 0.3703341175219492
$ws.Cells.Item(17, 12).Value = 0.8519797784390732
$ws.Cells.Item(17, 13).Value = 0.5819574805594669

$ws.Cells.Item(18, 2).Value = 1.605933551347164
$ws.Cells.Item(18, 4).Value = 0.178610402261171
$ws.Cells.Item(18, 5).Value = 0.2460700939134526
$ws.Cells.Item(18, 6).Value = 1.839408605385586
$ws.Cells.Item(18, 7).Value = 1.391667499417792
$ws.Cells.Item(18, 8).Value = 1.22162850036861
$ws.Cells.Item(18, 9).Value = 1.008356637970956
$ws.Cells.Item(18, 10).Value = 0.3653838029165399
$ws.Cells.Item(18, 12).Value = 0.837726038034134
$ws.Cells.Item(18, 13).Value = 0.5727656432507331

$ws.Cells.Item(19, 2).Value = 1.597657255757781
$ws.Cells.Item(19, 4).Value = 0.1779955408266574
$ws.Cells.Item(19, 5).Value = 0.2450387330885633
$ws.Cells.Item(19, 6).Value = 1.83940639040506
$ws.Cells.Item(19, 7).Value = 1.389989284778153
$ws.Cells.Item(19, 8).Value = 1.221443357717789
$ws.Cells.Item(19, 9).Value = 1.010067058448406
$ws.Cells.Item(19, 10).Value = 0.3637069242726483
$ws.Cells.Item(19, 12).Value = 0.8329015440551188
$ws.Cells.Item(19, 13).Value = 0.569654283656817

$ws.Cells.Item(20, 2).Value = 1.634908612726235
$ws.Cells.Item(20, 4).Value = 0.1807600154591853
$ws.Cells.Item(20, 5).Value = 0.2496774513424498
$ws.Cells.Item(20, 6).Value = 1.839518017158724
$ws.Cells.Item(20, 7).Value = 1.397630483917993
$ws.Cells.Item(20, 8).Value = 1.222337073006258
$ws.Cells.Item(20, 9).Value = 1.002423475069348
$ws.Cells.Item(20, 10).Value = 0.3712499339978876
$ws.Cells.Item(20, 12).Value = 0.8546185723551787
$ws.Cells.Item(20, 13).Value = 0.5836590789257201

$ws.Cells.Item(21, 2).Value = 1.760187473558631
$ws.Cells.Item(21, 4).Value = 0.1900053259945054
$ws.Cells.Item(21, 5).Value = 0.2652184477368706
$ws.Cells.Item(21, 6).Value = 1.841677947387367
$ws.Cells.Item(21, 7).Value = 1.42487020294476
$ws.Cells.Item(21, 8).Value = 1.226402621226867
$ws.Cells.Item(21, 9).Value = 0.9776797878672081
$ws.Cells.Item(21, 10).Value = 0.3965376517292896
$ws.Cells.Item(21, 12).Value = 0.9276929579073396
$ws.Cells.Item(21, 13).Value = 0.6307700401629575

$ws.Cells.Item(22, 2).Value = 1.842147532689978
$ws.Cells.Item(22, 4).Value = 0.1960160120913486
$ws.Cells.Item(22, 5).Value = 0.2753418257886082
$ws.Cells.Item(22, 6).Value = 1.844426979794676
$ws.Cells.Item(22, 7).Value = 1.443850090270246
$ws.Cells.Item(22, 8).Value = 1.229855587974157
$ws.Cells.Item(22, 9).Value = 0.96220903517926
$ws.Cells.Item(22, 10).Value = 0.4130219150217442
$ws.Cells.Item(22, 12).Value = 0.9755300586394924
$ws.Cells.Item(22, 13).Value = 0.661600271089128

$ws.Cells.Item(23, 2).Value = 1.798396679864197
$ws.Cells.Item(23, 4).Value = 0.1928108698676141
$ws.Cells.Item(23, 5).Value = 0.2699418859821492
$ws.Cells.Item(23, 6).Value = 1.84283694051507
$ws.Cells.Item(23, 7).Value = 1.433611829558203
$ws.Cells.Item(23, 8).Value = 1.227939605835275
$ws.Cells.Item(23, 9).Value = 0.9704018391018572
$ws.Cells.Item(23, 10).Value = 0.4042279362688816
$ws.Cells.Item(23, 12).Value = 0.9499914060361334
$ws.Cells.Item(23, 13).Value = 0.6451420063382756

$ws.Cells.Item(24, 2).Value = 1.632862481724715
$ws.Cells.Item(24, 4).Value = 0.1806083661648472
$ws.Cells.Item(24, 5).Value = 0.2494228812687922
$ws.Cells.Item(24, 6).Value = 1.839505190368939
$ws.Cells.Item(24, 7).Value = 1.397204997890242
$ws.Cells.Item(24, 8).Value = 1.222284005769751
$ws.Cells.Item(24, 9).Value = 1.002839700986534
$ws.Cells.Item(24, 10).Value = 0.370835914727877
$ws.Cells.Item(24, 12).Value = 0.853425565292838
$ws.Cells.Item(24, 13).Value = 0.5828897842285556

$ws.Cells.Item(25, 2).Value = 1.454701747099193
$ws.Cells.Item(25, 4).Value = 0.1673077788248492
$ws.Cells.Item(25, 5).Value = 0.2271484253553666
$ws.Cells.Item(25, 6).Value = 1.841607274228863
$ws.Cells.Item(25, 7).Value = 1.362925048198179
$ws.Cells.Item(25, 8).Value = 1.219576061878655
$ws.Cells.Item(25, 9).Value = 1.040825982114495
$ws.Cells.Item(25, 10).Value = 0.3346412626350457
$ws.Cells.Item(25, 12).Value = 0.7496135466936948
$ws.Cells.Item(25, 13).Value = 0.5159261788612284
